# feat: add 2022-Q1 data
#
# The workbook previously had quarterly detail sheets 2020-Q4 .. 2021-Q4
# plus a trailing "总计" (totals/history) summary sheet. This change adds
# a "2022-Q1" quarter:
#   * a new "2022-Q1" sheet (per-fund holdings detail, same shape as the
#     other quarterly sheets) inserted right before "总计"
#   * the "总计" sheet gets a new leading row summarising 2022-Q1
#
# To land the new sheet in the right slot with the right internal sheetId
# bookkeeping, the existing "总计" sheet is renamed to "2022-Q1" and
# repopulated with the quarter's fund-holdings detail, and a brand new
# "总计" sheet is appended at the end with the (updated) totals history.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Borrow styling (style index used for header row + index column) from an
# existing quarterly sheet via copy / paste-special so formats match
# exactly instead of synthesizing a near-duplicate style.
$donor = $wb.Worksheets.Item("2021-Q4")
$donor.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$donor.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B-G are text-typed in the source data (even the numeric-looking
# ones) - force text formatting before writing so "11.07" etc. are stored
# as text, matching the sibling quarter sheets. Reset the style back to
# "Normal" afterwards so the cells don't carry a leftover custom-format
# style index (the source sheets use the plain default style there).
$q1.Range("B2:G9").NumberFormat = "@"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "000893"
$q1.Range("C2").Value = "工银创新动力股票"
$q1.Range("D2").Value = "11.07"
$q1.Range("E2").Value = "81.01"
$q1.Range("F2").Value = "3.32"
$q1.Range("G2").Value = "0.3675"
$q1.Range("H2").Value = 7

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "519616"
$q1.Range("C3").Value = "银河君信灵活配置混合A"
$q1.Range("D3").Value = "4.54"
$q1.Range("E3").Value = "24.42"
$q1.Range("F3").Value = "0.98"
$q1.Range("G3").Value = "0.0445"
$q1.Range("H3").Value = 3

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "519618"
$q1.Range("C4").Value = "银河君信灵活配置混合I"
$q1.Range("D4").Value = "4.54"
$q1.Range("E4").Value = "24.42"
$q1.Range("F4").Value = "0.98"
$q1.Range("G4").Value = "0.0445"
$q1.Range("H4").Value = 3

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "003805"
$q1.Range("C5").Value = "华安新恒利灵活配置混合A"
$q1.Range("D5").Value = "5.58"
$q1.Range("E5").Value = "24.58"
$q1.Range("F5").Value = "0.44"
$q1.Range("G5").Value = "0.0246"
$q1.Range("H5").Value = 2

$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "519617"
$q1.Range("C6").Value = "银河君信灵活配置混合C"
$q1.Range("D6").Value = "0.64"
$q1.Range("E6").Value = "24.42"
$q1.Range("F6").Value = "0.98"
$q1.Range("G6").Value = "0.0063"
$q1.Range("H6").Value = 3

$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "003806"
$q1.Range("C7").Value = "华安新恒利灵活配置混合C"
$q1.Range("D7").Value = "1.20"
$q1.Range("E7").Value = "24.58"
$q1.Range("F7").Value = "0.44"
$q1.Range("G7").Value = "0.0053"
$q1.Range("H7").Value = 2

$q1.Range("A8").Value = 6
$q1.Range("B8").Value = "009387"
$q1.Range("C8").Value = "嘉实稳福混合A"
$q1.Range("D8").Value = "0.08"
$q1.Range("E8").Value = "34.71"
$q1.Range("F8").Value = "1.48"
$q1.Range("G8").Value = "0.0012"
$q1.Range("H8").Value = 10

$q1.Range("A9").Value = 7
$q1.Range("B9").Value = "009388"
$q1.Range("C9").Value = "嘉实稳福混合C"
$q1.Range("D9").Value = "0.01"
$q1.Range("E9").Value = "34.71"
$q1.Range("F9").Value = "1.48"
$q1.Range("G9").Value = "0.0001"
$q1.Range("H9").Value = 10

# Drop the transient "@" number-format style again, now that the text
# values are committed - restores the default (unstyled) cell style.
$q1.Range("B2:G9").Style = "Normal"

# ---------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet with the updated totals history
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

$donor.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$donor.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 0.49

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.51

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 8
$total.Range("D4").Value = 0.55

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 5
$total.Range("D5").Value = 0.15

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 11
$total.Range("D6").Value = 0.44

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 1
$total.Range("D7").Value = 0.05
